# The document carries three embedded pictures that live in the
# headers/footers (a Pearson logo repeated in the "default" and
# "first page" footers, and a BTEC logo in the "first page" header).
# Each one needs its display/shape name swapped:
#   - the two Pearson logo pictures: "image1.png" -> "image2.png"
#   - the BTEC logo picture:         "image2.jpg" -> "image1.jpg"
#
# Inline pictures that live in headers/footers are not reachable via
# Document.InlineShapes (that collection only covers the main story),
# so walk Sections(1).Headers/Footers and use each HeaderFooter's
# Range.InlineShapes collection instead. Shapes are matched by their
# alt text (the "descr" attribute), which is stable and read back
# correctly, rather than by .Name (write-only in this object model).

$d = $word.ActiveDocument
$section = $d.Sections(1)

$pearsonAlt = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"
$btecAlt = "BTec_Logo-Orange"

# --- Footers -------------------------------------------------------------
# The default and first-page footers each carry one inline Pearson
# Edexcel logo picture named "image1.png"; rename it to "image2.png".
for ($i = 1; $i -le 3; $i++) {
    $footer = $section.Footers($i)
    if ($footer.Exists) {
        for ($j = 1; $j -le $footer.Range.InlineShapes.Count; $j++) {
            $shape = $footer.Range.InlineShapes($j)
            if ($shape.AlternativeText -eq $pearsonAlt) {
                $shape.Name = "image2.png"
            }
        }
    }
}

# --- Headers ---------------------------------------------------------
# The first-page header carries the BTEC logo picture named
# "image2.jpg"; rename it to "image1.jpg".
for ($i = 1; $i -le 3; $i++) {
    $header = $section.Headers($i)
    if ($header.Exists) {
        for ($j = 1; $j -le $header.Range.InlineShapes.Count; $j++) {
            $shape = $header.Range.InlineShapes($j)
            if ($shape.AlternativeText -eq $btecAlt) {
                $shape.Name = "image1.jpg"
            }
        }
    }
}

Write-Output "Renamed header/footer logo pictures."
